# Generate Report for Handoff
#
# Re-runs the localization status report: the "Handed back" status becomes
# "Ready for handoff" everywhere it appears, the handoff timestamps advance
# to the new generation run, and the two handback files (one per language)
# are flagged as stale with an "Error Detail" message pointing at the
# latest handback commit on GitHub.

$wb = $excel.ActiveWorkbook

$statusNew = "Ready for handoff"
$genDate   = "2016-09-06 03:48:40"
$zhDate    = "2016-09-06 03:48:23"

$err4086 = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2f267e57b3bc0bffdd253d7bb7da685eedefb5a0/e2e/4086031d-1263-438b-b1ec-68f462519a3e.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d361491e552b207a1943654fdc92dad8051b421c/e2e/4086031d-1263-438b-b1ec-68f462519a3e.md."
$err8697 = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2f267e57b3bc0bffdd253d7bb7da685eedefb5a0/e2e/8697c96a-d033-4cb2-9976-0fb2cad765db.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d361491e552b207a1943654fdc92dad8051b421c/e2e/8697c96a-d033-4cb2-9976-0fb2cad765db.md."

$dateFormat = "yyyy-mm-dd HH:mm:ss"
$hyperlinkColor = 15570276   # BGR long for RGB(100,149,237) "Cornflower Blue", matches the sheet's HyperLink cell style

# --- Overview sheet -------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = $statusNew
$wsOverview.Range("F2").Value = $statusNew
$wsOverview.Range("G2").Value = $genDate

$wsOverview.Range("E3").Value = $statusNew
$wsOverview.Range("F3").Value = $statusNew
$wsOverview.Range("G3").Value = $genDate

# Columns E and F narrow from the old width to the new, tighter width.
$wsOverview.Range("E1").ColumnWidth = 16.3
$wsOverview.Range("F1").ColumnWidth = 16.3

# Preserve the existing look (date format / hyperlink styling) across the
# round-trip for the cells this edit doesn't otherwise touch.
foreach ($r in @("G2","G3")) {
    $wsOverview.Range($r).NumberFormat = $dateFormat
}
foreach ($r in @("B2","B3")) {
    $wsOverview.Range($r).Font.Underline = $true
    $wsOverview.Range($r).Font.Color = $hyperlinkColor
}

# --- zh-cn sheet ------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C2").Value = $statusNew
$wsZhCn.Range("H2").Value = $zhDate
$wsZhCn.Range("P2").Value = $err4086

$wsZhCn.Range("C3").Value = $statusNew
$wsZhCn.Range("H3").Value = $zhDate
$wsZhCn.Range("P3").Value = $err8697

$wsZhCn.Range("C1").ColumnWidth = 16.3
$wsZhCn.Range("P1").ColumnWidth = 39.166666666666664

foreach ($r in @("H2","K2","H3","K3")) {
    $wsZhCn.Range($r).NumberFormat = $dateFormat
}
foreach ($r in @("A2","I2","A3","I3")) {
    $wsZhCn.Range($r).Font.Underline = $true
    $wsZhCn.Range($r).Font.Color = $hyperlinkColor
}

# --- de-de sheet ------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C2").Value = $statusNew
$wsDeDe.Range("H2").Value = $genDate
$wsDeDe.Range("P2").Value = $err4086

$wsDeDe.Range("C3").Value = $statusNew
$wsDeDe.Range("H3").Value = $genDate
$wsDeDe.Range("P3").Value = $err8697

$wsDeDe.Range("C1").ColumnWidth = 16.3
$wsDeDe.Range("P1").ColumnWidth = 39.166666666666664

foreach ($r in @("H2","K2","H3","K3")) {
    $wsDeDe.Range($r).NumberFormat = $dateFormat
}
foreach ($r in @("A2","I2","A3","I3")) {
    $wsDeDe.Range($r).Font.Underline = $true
    $wsDeDe.Range($r).Font.Color = $hyperlinkColor
}
